$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes existing rows 3..9 down to 4..10)
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the new weekly record
$ws.Cells.Item(3, 1).Value = 8
$ws.Cells.Item(3, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(3, 3).Value = "Coquimbo"
$ws.Cells.Item(3, 4).Value = 44998
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(3, 6).Value = 100114002
$ws.Cells.Item(3, 7).Value = "Camote"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 320
$ws.Cells.Item(3, 11).Value = 17000
$ws.Cells.Item(3, 12).Value = 18000
$ws.Cells.Item(3, 13).Value = 17500
$ws.Cells.Item(3, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(3, 15).Value = "Perú"
$ws.Cells.Item(3, 16).Value = 972
$ws.Cells.Item(3, 17).Value = 18
$ws.Cells.Item(3, 18).Value = "Hortaliza"
